$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 4, shifting existing rows 4:75 down to 5:76
$ws.Rows.Item(4).Insert()

# Populate the newly inserted row 4 with the new record's data
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C4").Value = "Coquimbo"
$ws.Range("D4").Value = 44616
$ws.Range("E4").Value = 4
$ws.Range("F4").Value = 100112030
$ws.Range("G4").Value = "Poroto granado"
$ws.Range("H4").Value = "Sin especificar"
$ws.Range("I4").Value = "Primera"
$ws.Range("J4").Value = 600
$ws.Range("K4").Value = 21000
$ws.Range("L4").Value = 23000
$ws.Range("M4").Value = 22000
$ws.Range("N4").Value = "$/malla 25 kilos"
$ws.Range("O4").Value = "Provincia de Limarí"
$ws.Range("P4").Value = 880
$ws.Range("Q4").Value = 25
$ws.Range("R4").Value = "Hortaliza"
